$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Updated transition-probability matrix cells (St. Mary-Woods_A)
# Values reflect re-simulated games / updated win-probability matrix.

# Row 2
$ws.Range("B2").Value = 0.2941176470588235
$ws.Range("C2").Value = 0.3529411764705883
$ws.Range("P2").Value = 0.1764705882352941
$ws.Range("S2").Value = 0.1764705882352941

# Row 3
$ws.Range("C3").Value = 0.1428571428571428
$ws.Range("P3").Value = 0.4285714285714285
$ws.Range("S3").Value = 0.4285714285714285

# Row 6
$ws.Range("F6").Value = 0.09090909090909091
$ws.Range("J6").Value = 0.2727272727272727
$ws.Range("Q6").Value = 0.04545454545454546
$ws.Range("R6").Value = 0.04545454545454546
$ws.Range("S6").Value = 0.5454545454545454

# Row 7
$ws.Range("J7").Value = 0.125
$ws.Range("Q7").Value = 0.3125
$ws.Range("R7").Value = 0.0625
$ws.Range("S7").Value = 0.5

# Row 8
$ws.Range("B8").Value = 0.1282051282051282
$ws.Range("F8").Value = 0.07692307692307693
$ws.Range("J8").Value = 0.1538461538461539
$ws.Range("Q8").Value = 0.1025641025641026
$ws.Range("R8").Value = 0.1025641025641026
$ws.Range("S8").Value = 0.4358974358974359

# Row 9
$ws.Range("F9").Value = 0.04
$ws.Range("J9").Value = 0.24
$ws.Range("O9").Value = 0.08
$ws.Range("Q9").Value = 0.12
$ws.Range("R9").Value = 0.12
$ws.Range("S9").Value = 0.4

# Row 10
$ws.Range("B10").Value = 0.08139534883720931
$ws.Range("F10").Value = 0.09302325581395349
$ws.Range("J10").Value = 0.1046511627906977
$ws.Range("O10").Value = 0.02325581395348837
$ws.Range("Q10").Value = 0.1279069767441861
$ws.Range("R10").Value = 0.01162790697674419
$ws.Range("S10").Value = 0.5581395348837209

# Row 11
$ws.Range("G11").Value = 0.25
$ws.Range("J11").Value = 0.08333333333333333

# Row 12
$ws.Range("G12").Value = 0.75
$ws.Range("J12").Value = 0.25

# Row 15
$ws.Range("H15").Value = 0.2631578947368421
$ws.Range("I15").Value = 0.2105263157894737
$ws.Range("J15").Value = 0.2631578947368421
$ws.Range("S15").Value = 0.2631578947368421

# Row 16
$ws.Range("H16").Value = 0.3333333333333333
$ws.Range("I16").Value = 0.3333333333333333
$ws.Range("J16").Value = 0.1666666666666667
$ws.Range("O16").Value = 0.1666666666666667

# Row 17
$ws.Range("H17").Value = 0.08333333333333333
$ws.Range("I17").Value = 0.04166666666666666
$ws.Range("J17").Value = 0.2916666666666667
$ws.Range("K17").Value = 0.1666666666666667
$ws.Range("M17").Value = 0.04166666666666666
$ws.Range("O17").Value = 0.125
$ws.Range("S17").Value = 0.25

# Row 18
$ws.Range("H18").Value = 0.1
$ws.Range("I18").Value = 0.1
$ws.Range("J18").Value = 0.3
$ws.Range("O18").Value = 0.1
$ws.Range("S18").Value = 0.2

# Row 19
$ws.Range("F19").Value = 0.04545454545454546
$ws.Range("H19").Value = 0.2272727272727273
$ws.Range("I19").Value = 0.1287878787878788
$ws.Range("J19").Value = 0.2803030303030303
$ws.Range("K19").Value = 0.07575757575757576
$ws.Range("M19").Value = 0.02272727272727273
$ws.Range("O19").Value = 0.06060606060606061
$ws.Range("S19").Value = 0.1590909090909091
